{"js": "// Office.js (Word JavaScript API) implementation of the two text edits\n// described by the diff:\n//  1) \"...Menschen, da diese sich in den wenigsten...\" ->\n//     \"...Menschen, da sich diese in den wenigsten...\" (+ trailing space\n//     added at the very end of the sentence).\n//  2) The two paragraphs \"Da viele \u00e4ltere Menschen einsam leben, ...\" and\n//     \"Der Roboter sollte unserer Vorstellung nach ...\" are merged into a\n//     single, fully rewritten paragraph.\n\nconst body = context.document.body;\n\n// --- Change 1 -------------------------------------------------------\nconst oldSentence =\n  \"Auch moderne Online-Dienstleistungen entsprechen h\u00e4ufig nicht den \" +\n  \"Bed\u00fcrfnissen \u00e4lterer Menschen, da diese sich in den wenigsten F\u00e4llen \" +\n  \"ausreichend mit Computern auskennen.\";\nconst newSentence =\n  \"Auch moderne Online-Dienstleistungen entsprechen h\u00e4ufig nicht den \" +\n  \"Bed\u00fcrfnissen \u00e4lterer Menschen, da sich diese in den wenigsten F\u00e4llen \" +\n  \"ausreichend mit Computern auskennen. \";\n\nconst hits1 = body.search(oldSentence, { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\n\nif (hits1.items.length > 0) {\n  hits1.items[0].insertText(newSentence, \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2 -------------------------------------------------------\nconst hits2 = body.search(\"Da viele \u00e4ltere Menschen einsam leben,\", {\n  matchCase: true,\n});\nhits2.load(\"items\");\nawait context.sync();\n\nif (hits2.items.length > 0) {\n  const firstPara = hits2.items[0].paragraphs.getFirst();\n  const secondPara = firstPara.getNext();\n\n  const mergedText =\n    \"Um das Problem, welches das Tragen von Eink\u00e4ufen darstellt, zu \" +\n    \"l\u00f6sen, entschieden wir uns einen Roboter zu konstruieren welcher \" +\n    \"diese Aufgabe \u00fcbernehmen sollte. Dieser sollte, um f\u00fcr die meist \" +\n    \"Technik-Unerfahrene \u00e4ltere Generation leicht bedienbar zu sein, \" +\n    \"autonom, also ohne weitere manuelle Steuerung, seinem Besitzer \" +\n    \"folgen und ansonsten nur mit Sprachbefehlen gesteuert werden. \" +\n    \"Au\u00dferdem muss der Roboter sich nat\u00fcrlich den Gegebenheiten, wie \" +\n    \"zum Beispiel der u.U. geringen Geschwindigkeit \u00e4lterer Menschen \" +\n    \"anpassen k\u00f6nnen. \";\n\n  firstPara.getRange().insertText(mergedText, \"Replace\");\n  secondPara.delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop implementation of the two text edits described by the\n# diff:\n#  1) \"...Menschen, da diese sich in den wenigsten...\" ->\n#     \"...Menschen, da sich diese in den wenigsten...\" (+ trailing space\n#     added at the very end of the sentence).\n#  2) The two paragraphs \"Da viele \u00e4ltere Menschen einsam leben, ...\" and\n#     \"Der Roboter sollte unserer Vorstellung nach ...\" are merged into a\n#     single, fully rewritten paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 ---------------------------------------------------------\n$oldSentence = \"Auch moderne Online-Dienstleistungen entsprechen h\u00e4ufig nicht den Bed\u00fcrfnissen \u00e4lterer Menschen, da diese sich in den wenigsten F\u00e4llen ausreichend mit Computern auskennen.\"\n$newSentence = \"Auch moderne Online-Dienstleistungen entsprechen h\u00e4ufig nicht den Bed\u00fcrfnissen \u00e4lterer Menschen, da sich diese in den wenigsten F\u00e4llen ausreichend mit Computern auskennen. \"\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = $oldSentence\n$find1.Replacement.Text = $newSentence\n$find1.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null\n\n# --- Change 2 ---------------------------------------------------------\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \"Da viele \u00e4ltere Menschen einsam leben,\"\n$find2.Execute() | Out-Null\n\n# $rng2 itself collapses onto the matched text once Find.Execute succeeds,\n# so the paragraph that contains the match is simply its first paragraph.\n$firstPara = $rng2.Paragraphs(1)\n$mergedText = \"Um das Problem, welches das Tragen von Eink\u00e4ufen darstellt, zu l\u00f6sen, entschieden wir uns einen Roboter zu konstruieren welcher diese Aufgabe \u00fcbernehmen sollte. Dieser sollte, um f\u00fcr die meist Technik-Unerfahrene \u00e4ltere Generation leicht bedienbar zu sein, autonom, also ohne weitere manuelle Steuerung, seinem Besitzer folgen und ansonsten nur mit Sprachbefehlen gesteuert werden. Au\u00dferdem muss der Roboter sich nat\u00fcrlich den Gegebenheiten, wie zum Beispiel der u.U. geringen Geschwindigkeit \u00e4lterer Menschen anpassen k\u00f6nnen. \"\n$oldTail = \"Der Roboter sollte unserer Vorstellung nach mithilfe einer Transportfl\u00e4che beispielsweise Eink\u00e4ufe tragen k\u00f6nnen.\"\n\n# Clear the first paragraph's content (this merges it with the following\n# paragraph, since the paragraph mark disappears too) and insert the\n# replacement text at the same spot.\n$r = $firstPara.Range\n$r.Delete()\n$r.InsertBefore($mergedText)\n\n# The old \"Der Roboter ...\" text now trails right after $mergedText inside\n# the very same (merged) paragraph; locate and remove it via Find so the\n# paragraph ends up containing only $mergedText.\n$tailScope = $d.Range($firstPara.Range.Start, $firstPara.Range.End)\n$find3 = $tailScope.Find\n$find3.ClearFormatting()\n$find3.Text = $oldTail\nif ($find3.Execute($oldTail, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)) {\n  $tailScope.Delete()\n}\n"}
